$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The OOXML <col width> persisted by this host is ColumnWidth + 5/6, so to
# land on a stored width of exactly 34 we back out that offset here.
$ws.Columns.Item(1).ColumnWidth = 34 - (5/6)

$names = @{
    2  = "Honda_HR-V_2022_"
    3  = "Dacia_Jogger_2021_"
    4  = "DS_4_2022_Standard_"
    5  = "DS_4_2022_Safety_Pack_"
    6  = "Volvo_C40_Recharge_2022_"
    7  = "Mercedes-Benz_C-Class_2022_"
    8  = "Kia_EV6_2022_"
    9  = "VW_Multivan_2022_"
    10 = "Peugeot_308_2022_"
    11 = "ORA_Funky_Cat_2022_"
    12 = "WEY_Coffee_01_2022_"
    13 = "Kia_Niro_2022_Safety_Pack_"
    14 = "Kia_Niro_2022_Standard_"
    15 = "Genesis_GV60_2022_"
    16 = "Tesla_Model_Y_2022_"
    17 = "Land_Rover_Discovery_Sport_2022_"
    18 = "CHERY_OMODA5_2022_"
    19 = "VW_ID_Buzz_2022_"
    20 = "VW_Touran_2022_"
    21 = "Lexus_RX_2022_"
    22 = "Škoda_Octavia_2022_"
    23 = "MG_4_Electric_2022_"
    24 = "Ford_Ranger_2022_"
    25 = "Lucid_Air_2022_"
    26 = "Jeep_Grand Cherokee_2022_"
    27 = "VW_Amarok_2022_"
    28 = "Mercedes-Benz_GLC_2022_"
    29 = "Maxus_MIFA_9_2022_"
    30 = "Ford_Puma_2022_"
    31 = "BMW_2_Series_Coupe_2022_"
    32 = "Renault_Megane_E-Tech_2022_"
    33 = "VW_Polo_2022_"
    34 = "Lexus_NX_2022_"
    35 = "VW_Taigo_2022_"
    36 = "VW_Polo_2022_"
    37 = "Nissan_Ariya_2022_"
    38 = "smart_1_2022_"
    39 = "Hyundai_IONIQ_6_2022_"
    40 = "Toyota_Corolla Cross_2022_"
    41 = "Range_Rover_Sport_2022_"
    42 = "Isuzu_D-MAX_Crew_Cab_2022_"
    43 = "NIO_ET7_2022_"
    44 = "Range_Rover_2022_"
    45 = "Renault_Austral_2022_"
    46 = "DS_9_2022_"
    47 = "Tesla_Model_S_2022_"
    48 = "Honda_Civic_2022_"
    49 = "Nissan_X_Trail_2021_"
    50 = "WEY_Coffee_02_2022_"
    51 = "Toyota_bZ4X_2022_"
    52 = "SEAT_Ibiza_2022_"
    53 = "BMW_X1_2022_"
    54 = "Mobilize_Limo_2022_"
    55 = "Mercedes-EQ_EQE_2022_"
    56 = "BYD_ATTO_3_2022_"
    57 = "Citroen_C5_X_2022_"
    58 = "SEAT_Arona_2022_"
    59 = "MAZDA_CX-60_2022_"
    60 = "BMW_2_Series_Active_Tourer_2022_"
    61 = "VW_Golf_2022_"
    62 = "Kia_Sportage_2022_"
    63 = "BMW_i4_2022_"
    64 = "Mercedes-Benz_T-Class_2022_"
    65 = "Toyota_Aygo_X_2022_"
    66 = "Alfa_Romeo_Tonale_2022_"
    67 = "Cupra_Born_2022_"
}

foreach ($row in $names.Keys) {
    $ws.Cells.Item($row, 1).Value = $names[$row]
}
